# TaskJuggle.xlsx — "new try to order tasks"
# 1. Rename Sheet1 -> "Tasks juggle"
# 2. Add a new sheet "Repeated Tasks" after it
# 3. Add an "I" column on the main sheet listing a few repeating tasks
#    (mirrored onto the new "Repeated Tasks" sheet)
# 4. Tweak a couple of row heights + the remembered selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the original sheet ---------------------------------------
$ws.Name = "Tasks juggle"

# --- add the new sheet right after it ---------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Repeated Tasks"

# --- new entries in the "I" column of the main sheet -------------------
# Copy the formatting already used by I4/I5 (wrapped, 8pt Arial) onto the
# new cells before writing their text.
$ws.Range("I4").Copy()
$ws.Range("I8:I10").PasteSpecial(-4122)

$ws.Range("I8").Value = "change log.php"
$ws.Range("I9").Value = "Doc / Clean code "
$ws.Range("I10").Value = "User Doc"

# I11 keeps the plain/default formatting (not the wrapped style used above)
$ws.Range("I11").Font.Name = "Arial"
$ws.Range("I11").Font.Size = 10
$ws.Range("I11").Font.Bold = $false
$ws.Range("I11").WrapText = $false
$ws.Range("I11").VerticalAlignment = -4107
$ws.Range("I11").HorizontalAlignment = 1
$ws.Range("I11").Value = "Release →…"

# --- row height tweaks (rows 8 & 9 now match row 10's 28.35) -----------
$ws.Rows.Item(8).RowHeight = 28.35
$ws.Rows.Item(9).RowHeight = 28.35

# --- mirror the three repeating tasks onto the new sheet ---------------
$ws.Range("I4").Copy()
$ws2.Range("B3:B5").PasteSpecial(-4122)

$ws2.Range("B3").Value = "change log.php"
$ws2.Range("B4").Value = "Doc / Clean code "
$ws2.Range("B5").Value = "User Doc"

$ws2.Rows.Item(4).RowHeight = 19.5

$ws2.Range("B3").Select() | Out-Null

# --- restore focus/selection on the main sheet --------------------------
$ws.Activate() | Out-Null
$ws.Range("I12").Select() | Out-Null
